$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.997.56"
$ws.Range("E2").Value = "  +0.12%  "

$ws.Range("D3").Value = "2.585.65"
$ws.Range("E3").Value = "  +1.72%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.42"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.63"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.34%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  +2.52%  "

$ws.Range("E9").Value = "  +2.72%  "

$ws.Range("E10").Value = "  +2.97%  "

$ws.Range("E11").Value = "  +0.06%  "

$ws.Range("E12").Value = "  -0.17%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.30"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.11%  "

$ws.Range("D14").Value = "3.049.91"
$ws.Range("E14").Value = "  +1.83%  "

$ws.Range("D15").Value = "62.916.20"
$ws.Range("E15").Value = "  +0.11%  "

$ws.Range("E16").Value = "  +3.34%  "

$ws.Range("D17").Value = "2.589.05"
$ws.Range("E17").Value = "  +2.72%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.32"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "343.16"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.50%  "

$ws.Range("E20").Value = "  +1.58%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.66"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.00%  "

$ws.Range("E22").Value = "  -0.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.68"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.31%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.33"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.36%  "

$ws.Range("D25").Value = "2.722.01"
$ws.Range("E25").Value = "  +2.33%  "

$ws.Range("E26").Value = "  -0.74%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.60"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.60%  "

$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.33"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.84"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +8.37%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.43"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.28%  "

$ws.Range("E32").Value = "  +4.09%  "

$ws.Range("D33").Value = "0.0₃0824"
$ws.Range("E33").Value = "  +1.65%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "467.02"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +16.83%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "175.01"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.50%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.60"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +4.58%  "

$ws.Range("E37").Value = "  +0.02%  "

$ws.Range("E38").Value = "  +1.95%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.13"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.50%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.55"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +5.44%  "

$ws.Range("E42").Value = "  -1.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "158.42"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.33%  "

$ws.Range("E44").Value = "  +1.40%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.638"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +6.69%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.25"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.69%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0544"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0969"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.90%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0237"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.56"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.51%  "

$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.43"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.04%  "
